# Update column F ("dSF") values for the rows that were repulled/recomputed.
# Column F previously mirrored column E ("dS0"); this commit pushes in the
# freshly re-pulled dSF values (and related mean-calculation updates).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    4  = -8
    5  = -6
    6  = 3
    12 = -5
    13 = -3
    15 = -6
    20 = 1
    21 = -6
    26 = -8
    27 = -3
    29 = -4
    30 = 0
    32 = -6
    33 = -4
    34 = 7
    40 = -6
    47 = -6
    48 = -5
    54 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
